$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename rain model constants (shared string values in column A, rows 19-20)
$ws.Range("A20").Value = "fi_lidar_rain_intensity"
$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"

# Update selection to match the new active cell
$ws.Range("A20").Select()
